$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.551.10'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '1.956.28'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '243.97'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('E6').Value = '  -0.96%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '58.58'
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +2.42%  '
$ws.Range('E10').Value = '  -4.38%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '22.08'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').Value = '2.243.21'
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.824'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '13.63'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.26'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').Value = '1.958.75'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').Value = '36.424.99'
$ws.Range('E18').Value = '  -0.09%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.70'
$ws.Range('E19').Value = '  -0.41%  '
$ws.Range('D20').Value = '0.0₃0855'
$ws.Range('E20').Value = '  -1.15%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '228.20'
$ws.Range('E21').Value = '  -0.84%  '
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('E25').Value = '  +1.63%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.27'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('E27').Value = '  -1.01%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '160.12'
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.38'
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('E32').Value = '  -0.48%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0617'
$ws.Range('E33').Value = '  -2.62%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.28'
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.24'
$ws.Range('E36').Value = '  +3.53%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.33'
$ws.Range('E37').Value = '  +8.70%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.78'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('E39').Value = '  -8.34%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0980'
$ws.Range('E40').Value = '  -0.80%  '
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('E42').Value = '  -1.43%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0211'
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '15.94'
$ws.Range('E44').Value = '  -1.24%  '
$ws.Range('D45').Value = '1.361.34'
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('E46').Value = '  -1.07%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '87.64'
$ws.Range('E47').Value = '  -0.91%  '
$ws.Range('E48').Value = '  -1.36%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.83'
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('D50').Value = '2.133.73'
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '43.57'
$ws.Range('E51').Value = '  -5.35%  '
